$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# The "Periodo Mora" (E) table for trabajador ISIDORO SALAS AGUIRRE (CC 15609536)
# is being re-sorted into ascending period order (1804..1807), and a new
# record for PATRICIA MARTINEZ MUTIS (CC 45488904, period 1804) is inserted
# right after the first ISIDORO 1804 row -- this is "parte 1 de nuevos
# estado de cuenta" from the commit message: new EC records added and the
# worker/period rows re-ordered/updated.

# Row 16: ISIDORO SALAS AGUIRRE, period 1804, Valor Mora 69600 (was period 1807 @ 72000)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "15609536"
$ws.Range("D16").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E16").Value = "1804"
$ws.Range("F16").Value = 69600
$ws.Range("G16").Value = 1800000

# Row 17: new record -- PATRICIA MARTINEZ MUTIS, period 1804
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45488904"
$ws.Range("D17").Value = "PATRICIA MARTINEZ MUTIS"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 17709
$ws.Range("G17").Value = 781300

# Row 18: ISIDORO SALAS AGUIRRE, period 1805 (was period 1806)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "15609536"
$ws.Range("D18").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 72000
$ws.Range("G18").Value = 1800000

# Row 19: ISIDORO SALAS AGUIRRE, period 1806 (was period 1805)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "15609536"
$ws.Range("D19").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 72000
$ws.Range("G19").Value = 1800000

# Row 20: ISIDORO SALAS AGUIRRE, period 1807 (was period 1804 for PATRICIA @ 17709/781300)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "15609536"
$ws.Range("D20").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E20").Value = "1807"
$ws.Range("F20").Value = 72000
$ws.Range("G20").Value = 1800000

$wb.Save()
